{"js": "// Apply targeted text replacements: update the date paragraph and all\n// arithmetic expressions in the table, matching the exact old -> new pairs\n// from the source diff. Uses Body.search for each distinct (unique) old\n// string and InsertLocation.replace so run/paragraph formatting (fonts,\n// size, alignment) is preserved.\nconst replacements = [\n  [\"2023-01-24 Tuesday\", \"2023-01-25 Wednesday\"],\n  [\"6+66=\", \"68+16=\"],\n  [\"66-53=\", \"85-50=\"],\n  [\"99-28=\", \"73-60=\"],\n  [\"34-0=\", \"34+35=\"],\n  [\"18+58=\", \"99-66=\"],\n  [\"5+34=\", \"80-45=\"],\n  [\"35+28=\", \"1+49=\"],\n  [\"86-27=\", \"99-24=\"],\n  [\"48+2=\", \"58-46=\"],\n  [\"57+2=\", \"56-43=\"],\n  [\"55+12=\", \"85-30=\"],\n  [\"3+29=\", \"91-64=\"],\n  [\"81-75=\", \"38-14=\"],\n  [\"63-42=\", \"29+51=\"],\n  [\"96-75=\", \"36-19=\"],\n  [\"4+86=\", \"63+14=\"],\n  [\"64+20=\", \"60-57=\"],\n  [\"13+77=\", \"95-25=\"],\n  [\"88-53=\", \"31-1=\"],\n  [\"68+31=\", \"14+34=\"],\n  [\"42+30=\", \"34-32=\"],\n  [\"64+9=\", \"63-36=\"],\n  [\"71-18=\", \"50+1=\"],\n  [\"72-72=\", \"58-41=\"],\n  [\"86-45=\", \"28+26=\"],\n  [\"3+74=\", \"86-71=\"],\n  [\"17-12=\", \"87+4=\"],\n  [\"82+4=\", \"81-69=\"],\n  [\"48-23=\", \"43+48=\"],\n  [\"14-4=\", \"5+50=\"],\n  [\"96-69=\", \"6+17=\"],\n  [\"34+10=\", \"70+22=\"],\n  [\"16+0=\", \"92-71=\"],\n  [\"87-32=\", \"83-62=\"],\n  [\"48+46=\", \"59+7=\"],\n  [\"72+9=\", \"86-34=\"],\n  [\"0+6=\", \"22-16=\"],\n  [\"18-14=\", \"75-7=\"],\n  [\"54-41=\", \"56+24=\"],\n  [\"65-45=\", \"93-55=\"],\n  [\"36+23=\", \"55-26=\"],\n  [\"82+15=\", \"39+41=\"],\n  [\"76-37=\", \"82-64=\"],\n  [\"99-29=\", \"95+3=\"],\n  [\"61+16=\", \"24+22=\"],\n  [\"10+38=\", \"15+76=\"],\n  [\"15+65=\", \"11+74=\"],\n  [\"39-26=\", \"62-45=\"],\n  [\"43+3=\", \"2+83=\"],\n  [\"91-27=\", \"7+32=\"],\n  [\"70-4=\", \"26+50=\"],\n  [\"41+11=\", \"63-31=\"],\n  [\"76-62=\", \"40+54=\"],\n  [\"2+45=\", \"53-19=\"],\n  [\"2+43=\", \"44-9=\"],\n  [\"57-55=\", \"49-34=\"],\n  [\"40+38=\", \"51-37=\"],\n  [\"85+11=\", \"73-21=\"],\n  [\"1+22=\", \"99-96=\"],\n  [\"43-25=\", \"44-10=\"],\n  [\"10+28=\", \"71-44=\"],\n  [\"83-44=\", \"96-59=\"],\n  [\"82-19=\", \"51+14=\"],\n  [\"24+17=\", \"76+19=\"],\n  [\"10+1=\", \"84-40=\"],\n  [\"47+43=\", \"58+8=\"],\n  [\"75-9=\", \"53-16=\"],\n  [\"46+9=\", \"8+36=\"],\n  [\"1+97=\", \"81-37=\"],\n  [\"24+24=\", \"31-9=\"],\n  [\"5+64=\", \"25+68=\"],\n  [\"2+90=\", \"87-29=\"],\n  [\"11+82=\", \"49-39=\"],\n  [\"94-16=\", \"99-18=\"],\n  [\"26+22=\", \"24+18=\"],\n  [\"51-20=\", \"10+44=\"],\n  [\"93-87=\", \"11+50=\"],\n  [\"57-17=\", \"0+16=\"],\n  [\"51+12=\", \"11-3=\"],\n  [\"96-55=\", \"50+27=\"],\n  [\"31+41=\", \"30+48=\"],\n  [\"97-7=\", \"47-45=\"],\n  [\"30-14=\", \"8+44=\"],\n  [\"30+2=\", \"0+56=\"],\n  [\"49+11=\", \"81+17=\"],\n  [\"36-34=\", \"41+31=\"],\n  [\"95-85=\", \"79-27=\"],\n  [\"62+4=\", \"80-3=\"],\n  [\"51-28=\", \"70+4=\"],\n  [\"32-14=\", \"32+18=\"],\n  [\"34-33=\", \"71-8=\"],\n  [\"28+61=\", \"43+31=\"],\n  [\"88-29=\", \"84+0=\"],\n  [\"81-65=\", \"88-84=\"],\n  [\"82-29=\", \"85-17=\"],\n  [\"18+79=\", \"77-12=\"],\n  [\"24+48=\", \"41-29=\"],\n  [\"81-57=\", \"40+15=\"],\n  [\"2+61=\", \"39+60=\"],\n  [\"98-59=\", \"83-33=\"],\n];\n\nconst pending = [];\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n    matchWildcards: false\n  });\n  results.load(\"items\");\n  pending.push({ results, newText, oldText });\n}\nawait context.sync();\n\nfor (const { results, newText, oldText } of pending) {\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Apply targeted text replacements: update the date paragraph and all\n# arithmetic expressions in the table, matching the exact old -> new pairs\n# from the source diff. Uses Find/Replace on $d.Content for each distinct\n# (unique) old string so existing run/paragraph formatting (fonts, size,\n# alignment) is preserved.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-01-24 Tuesday\", \"2023-01-25 Wednesday\"),\n    @(\"6+66=\", \"68+16=\"),\n    @(\"66-53=\", \"85-50=\"),\n    @(\"99-28=\", \"73-60=\"),\n    @(\"34-0=\", \"34+35=\"),\n    @(\"18+58=\", \"99-66=\"),\n    @(\"5+34=\", \"80-45=\"),\n    @(\"35+28=\", \"1+49=\"),\n    @(\"86-27=\", \"99-24=\"),\n    @(\"48+2=\", \"58-46=\"),\n    @(\"57+2=\", \"56-43=\"),\n    @(\"55+12=\", \"85-30=\"),\n    @(\"3+29=\", \"91-64=\"),\n    @(\"81-75=\", \"38-14=\"),\n    @(\"63-42=\", \"29+51=\"),\n    @(\"96-75=\", \"36-19=\"),\n    @(\"4+86=\", \"63+14=\"),\n    @(\"64+20=\", \"60-57=\"),\n    @(\"13+77=\", \"95-25=\"),\n    @(\"88-53=\", \"31-1=\"),\n    @(\"68+31=\", \"14+34=\"),\n    @(\"42+30=\", \"34-32=\"),\n    @(\"64+9=\", \"63-36=\"),\n    @(\"71-18=\", \"50+1=\"),\n    @(\"72-72=\", \"58-41=\"),\n    @(\"86-45=\", \"28+26=\"),\n    @(\"3+74=\", \"86-71=\"),\n    @(\"17-12=\", \"87+4=\"),\n    @(\"82+4=\", \"81-69=\"),\n    @(\"48-23=\", \"43+48=\"),\n    @(\"14-4=\", \"5+50=\"),\n    @(\"96-69=\", \"6+17=\"),\n    @(\"34+10=\", \"70+22=\"),\n    @(\"16+0=\", \"92-71=\"),\n    @(\"87-32=\", \"83-62=\"),\n    @(\"48+46=\", \"59+7=\"),\n    @(\"72+9=\", \"86-34=\"),\n    @(\"0+6=\", \"22-16=\"),\n    @(\"18-14=\", \"75-7=\"),\n    @(\"54-41=\", \"56+24=\"),\n    @(\"65-45=\", \"93-55=\"),\n    @(\"36+23=\", \"55-26=\"),\n    @(\"82+15=\", \"39+41=\"),\n    @(\"76-37=\", \"82-64=\"),\n    @(\"99-29=\", \"95+3=\"),\n    @(\"61+16=\", \"24+22=\"),\n    @(\"10+38=\", \"15+76=\"),\n    @(\"15+65=\", \"11+74=\"),\n    @(\"39-26=\", \"62-45=\"),\n    @(\"43+3=\", \"2+83=\"),\n    @(\"91-27=\", \"7+32=\"),\n    @(\"70-4=\", \"26+50=\"),\n    @(\"41+11=\", \"63-31=\"),\n    @(\"76-62=\", \"40+54=\"),\n    @(\"2+45=\", \"53-19=\"),\n    @(\"2+43=\", \"44-9=\"),\n    @(\"57-55=\", \"49-34=\"),\n    @(\"40+38=\", \"51-37=\"),\n    @(\"85+11=\", \"73-21=\"),\n    @(\"1+22=\", \"99-96=\"),\n    @(\"43-25=\", \"44-10=\"),\n    @(\"10+28=\", \"71-44=\"),\n    @(\"83-44=\", \"96-59=\"),\n    @(\"82-19=\", \"51+14=\"),\n    @(\"24+17=\", \"76+19=\"),\n    @(\"10+1=\", \"84-40=\"),\n    @(\"47+43=\", \"58+8=\"),\n    @(\"75-9=\", \"53-16=\"),\n    @(\"46+9=\", \"8+36=\"),\n    @(\"1+97=\", \"81-37=\"),\n    @(\"24+24=\", \"31-9=\"),\n    @(\"5+64=\", \"25+68=\"),\n    @(\"2+90=\", \"87-29=\"),\n    @(\"11+82=\", \"49-39=\"),\n    @(\"94-16=\", \"99-18=\"),\n    @(\"26+22=\", \"24+18=\"),\n    @(\"51-20=\", \"10+44=\"),\n    @(\"93-87=\", \"11+50=\"),\n    @(\"57-17=\", \"0+16=\"),\n    @(\"51+12=\", \"11-3=\"),\n    @(\"96-55=\", \"50+27=\"),\n    @(\"31+41=\", \"30+48=\"),\n    @(\"97-7=\", \"47-45=\"),\n    @(\"30-14=\", \"8+44=\"),\n    @(\"30+2=\", \"0+56=\"),\n    @(\"49+11=\", \"81+17=\"),\n    @(\"36-34=\", \"41+31=\"),\n    @(\"95-85=\", \"79-27=\"),\n    @(\"62+4=\", \"80-3=\"),\n    @(\"51-28=\", \"70+4=\"),\n    @(\"32-14=\", \"32+18=\"),\n    @(\"34-33=\", \"71-8=\"),\n    @(\"28+61=\", \"43+31=\"),\n    @(\"88-29=\", \"84+0=\"),\n    @(\"81-65=\", \"88-84=\"),\n    @(\"82-29=\", \"85-17=\"),\n    @(\"18+79=\", \"77-12=\"),\n    @(\"24+48=\", \"41-29=\"),\n    @(\"81-57=\", \"40+15=\"),\n    @(\"2+61=\", \"39+60=\"),\n    @(\"98-59=\", \"83-33=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"WARNING: pattern not found: $oldText\"\n    }\n}\n"}
